$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("01_IB전략컨설팅부")
$ws.Rows("20:21").Delete()
